$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title (paragraph 1, single run) - simple text swap
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The Nexus of Art and Mathematics", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Exploring the Wonders of Chemistry: Unraveling the Secrets of Matter", 2
) | Out-Null

# ---------------------------------------------------------------------
# 2. Author paragraph (paragraph 2): "Emily Carter" becomes three runs:
#    "Dr" + "." + " Jane Carter" (same rPr). Rebuild the paragraph's run
#    content via InsertXML (keeps the paragraph's pPr / NoSpacing style).
# ---------------------------------------------------------------------
$authorRunsXml = @'
<w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t>Dr</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t xml:space="preserve"> Jane Carter</w:t>
      </w:r>
    
'@

$p2Full = $d.Paragraphs.Item(2).Range
$p2Target = $d.Range($p2Full.Start, $p2Full.End - 1)
$p2Target.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $authorRunsXml + '</w:p>')

# ---------------------------------------------------------------------
# 3. Email paragraph (paragraph 3): two simple text swaps
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "emilycarter60@domainhost", $true, $false, $false, $false, $false,
    $true, 1, $false, "janecarter09@educonnect", 2
) | Out-Null

$d.Content.Find.Execute(
    "org", $true, $true, $false, $false, $false,
    $true, 1, $false, "com", 2
) | Out-Null

# ---------------------------------------------------------------------
# 4. Body paragraph (paragraph 5): rewritten prose. Several sentences
#    swap 1:1 inside their existing run, two sentences that used to be
#    three runs collapse into one run, and two spots gain two brand
#    new runs (a "." run plus a following sentence run). Because the
#    run layout changes in the middle of this long paragraph, the most
#    reliable approach is to rebuild the whole paragraph's run content
#    in one shot via InsertXML (this preserves the paragraph's pPr,
#    which here is empty/default).
# ---------------------------------------------------------------------
$bodyRunsXml = @'
<w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>Chemistry, the study of matter and its properties, is an intriguing field that holds the key to understanding the world around us</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> It is a dynamic discipline that constantly evolves, offering new insights into the nature of matter and its interactions</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> In this essay, we will delve into the fascinating world of chemistry, uncovering the secrets behind the composition, structure, and behavior of matter</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>As we embark on this journey, we will explore the fundamental principles that govern chemical reactions, delving into the intricate dance of atoms and molecules as they rearrange to form new substances</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> We will uncover the secrets of chemical bonding, investigating the forces that hold atoms together and determine the properties of materials</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Moreover, we will unravel the mysteries of chemical reactions, examining the conditions under which they occur and the factors that influence their rates and outcomes</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Our adventure into the realm of chemistry will take us beyond the classroom, as we explore the practical applications of this science in everyday life</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> We will investigate the role of chemistry in the development of new materials, unraveling the processes behind the creation of plastics, metals, and ceramics</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> We will delve into the world of pharmaceuticals, examining the intricate mechanisms by which drugs interact with the human body to combat diseases</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Furthermore, we will explore the impact of chemistry on agriculture, uncovering the secrets behind fertilizers and pesticides that help to feed a growing population</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    
'@

$p5Full = $d.Paragraphs.Item(5).Range
$p5Target = $d.Range($p5Full.Start, $p5Full.End - 1)
$p5Target.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyRunsXml + '</w:p>')

# ---------------------------------------------------------------------
# 5. Summary paragraph (paragraph 7): rewritten prose, same pattern as
#    above - one spot gains a lastRenderedPageBreak run split, one spot
#    collapses two sentences + a "." run into a single run.
# ---------------------------------------------------------------------
$summaryRunsXml = @'
<w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>Chemistry, the study of matter and its properties, offers a fascinating exploration into the nature of the world around us</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> By unraveling the secrets of chemical reactions, investigating the principles of chemical bonding, and understanding the practical applications of chemistry in various fields, we gain a deeper appreciation for the intricate symphony of matter that governs </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>our existence</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> The study of chemistry not only enhances our knowledge of the material world but also equips us with the tools to address global challenges and shape a sustainable future</w:t>
      </w:r>
      <w:r w:rsidR="00AA66D6">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    
'@

$p7Full = $d.Paragraphs.Item(7).Range
$p7Target = $d.Range($p7Full.Start, $p7Full.End - 1)
$p7Target.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $summaryRunsXml + '</w:p>')

# ---------------------------------------------------------------------
# 6. A new, fully empty paragraph is appended at the very end of the
#    body (after the Summary paragraph, before the sectPr).
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

Write-Output "Stage A done"
Write-Output $d.Paragraphs.Item(1).Range.Text
Write-Output $d.Paragraphs.Item(2).Range.Text
Write-Output $d.Paragraphs.Item(3).Range.Text
Write-Output $d.Paragraphs.Item(5).Range.Text
Write-Output $d.Paragraphs.Item(7).Range.Text
Write-Output "Paragraphs.Count=$($d.Paragraphs.Count)"
